# Printlab Beverage Can Opener BOM - add GitHub link
# Replaces the "Add GitHub Link" placeholder in G7 with the actual
# GitHub URL and turns it into a real hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$githubUrl = "https://github.com/makersmakingchange/Printlab_Beverage_Can_Opener/tree/main/Build_Files/3D_Print_Files"

$linkCell = $ws.Range("G7")

# Remember the existing formatting (the cell was already styled like a
# hyperlink as a placeholder) so it is preserved after Hyperlinks.Add
# applies its own style.
$origStyle = $linkCell.Style

$linkCell.Value2 = $githubUrl

$ws.Hyperlinks.Add($linkCell, $githubUrl) | Out-Null

$linkCell.Style = $origStyle

# Move the active selection, matching where the author ended up after
# making the edit.
$ws.Range("H16").Select() | Out-Null

$wb.Save() | Out-Null
